$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# B4 (2018 row): 23 (2017) + 40 -> 63 (was +38 -> 61)
$ws.Range("B4").Formula = "=B3+40"

# B5 (2019 row): 63 (2018) + 55 -> 118 (was +54 -> 115)
$ws.Range("B5").Formula = "=B4+55"

# Leave the cursor on B5, matching the saved selection in the edited file
$ws.Range("B5").Select()
